$wb = $excel.ActiveWorkbook

# --- Rename sheets ---
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws1.Name = "Field Breakdown"
$ws2.Name = "ID Lookup"

# --- Add new columns F:H on "Field Breakdown" (mirrors A:C with lowercase gender labels) ---
$ws1.Range("F1").Value = "Disease_Cases"
$ws1.Range("G1").Value = "Census_Population"
$ws1.Range("H1").Value = "Projected_Population"
$ws1.Range("F1:H1").Font.Bold = $true

$ws1.Range("F2").Value = "All"
$ws1.Range("G2").Value = "both sexes"
$ws1.Range("H2").Value = "both sexes"

$ws1.Range("F3").Value = "Female"
$ws1.Range("G3").Value = "female"
$ws1.Range("H3").Value = "female"

$ws1.Range("F4").Value = "Male"
$ws1.Range("G4").Value = "male"
$ws1.Range("H4").Value = "male"

# Empty bold cells below the mini header block (matches the blank divider row under A:C)
$ws1.Range("F14:H14").Font.Bold = $true

# Approximate the auto-fit column widths for the new columns
$ws1.Columns.Item(6).ColumnWidth = 12.0
$ws1.Columns.Item(7).ColumnWidth = 16.5
$ws1.Columns.Item(8).ColumnWidth = 18.65

# --- Sheet view / selection changes ---
# Set "ID Lookup" selection first (it is no longer the active tab, so it must
# not be the last sheet activated/selected).
$ws2.Activate() | Out-Null
$ws2.Range("B30").Select() | Out-Null

# "Field Breakdown" becomes the active sheet/tab, with D10 selected. Doing
# this last ensures it ends up as the active sheet in the saved workbook.
$ws1.Activate() | Out-Null
$ws1.Range("D10").Select() | Out-Null
